# AutoCommit_5 октября 2023 г. 19:04:55_SibNout2023
#
# The "Date Placeholder" shapes on the slide master and every slide
# layout contain a cached <a:fld type="datetimeFigureOut"> whose display
# text was last refreshed on 06.02.2023. Re-stamp every one of them with
# the date this automated commit actually happened on: 05.10.2023
# (05 October 2023), matching PowerPoint's "update automatically" date
# placeholder being re-rendered/re-saved.

$p = $ppt.ActivePresentation
$newDate = "05.10.2023"

# Map of CustomLayout index -> name of its date placeholder shape.
$dateShapeNameByLayout = @{
    1  = "Date Placeholder 3"
    2  = "Date Placeholder 3"
    3  = "Date Placeholder 3"
    4  = "Date Placeholder 4"
    5  = "Date Placeholder 6"
    6  = "Date Placeholder 2"
    7  = "Date Placeholder 1"
    8  = "Date Placeholder 4"
    9  = "Date Placeholder 4"
    10 = "Date Placeholder 4"
    11 = "Date Placeholder 3"
    12 = "Date Placeholder 3"
    13 = "Date Placeholder 3"
    14 = "Date Placeholder 6"
    15 = "Date Placeholder 6"
    16 = "Date Placeholder 3"
    17 = "Date Placeholder 3"
}

$master = $p.SlideMaster

# Slide master's own date placeholder.
$masterDateShape = $master.Shapes.Item("Date Placeholder 3")
$masterDateShape.TextFrame.TextRange.Text = $newDate

# Every custom (slide) layout's date placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapeName = $dateShapeNameByLayout[$i]
    $shp = $layout.Shapes.Item($shapeName)
    $shp.TextFrame.TextRange.Text = $newDate
}
